$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Wyniki")

# Row 2 - RandomForestClassifier
$ws.Range("C2").Value = 0.845
$ws.Range("D2").Value = 0.8625
$ws.Range("E2").Value = 0.89
$ws.Range("F2").Value = 0.8825
$ws.Range("G2").Value = 0.88
$ws.Range("H2").Value = 0.855
$ws.Range("I2").Value = 0.8739999999999999

# Row 3 - KNeighborsClassifier
$ws.Range("C3").Value = 0.78

# Row 4 - MultinomialNB
$ws.Range("C4").Value = 0.745

# Row 5 - DecisionTreeClassifier
$ws.Range("B5").Value = 0.845
$ws.Range("C5").Value = 0.7524999999999999
$ws.Range("D5").Value = 0.8
$ws.Range("E5").Value = 0.8149999999999999
$ws.Range("F5").Value = 0.8
$ws.Range("G5").Value = 0.8275
$ws.Range("H5").Value = 0.7675
$ws.Range("I5").Value = 0.8019999999999999

# Row 6 - BernoulliNB
$ws.Range("C6").Value = 0.6925

# Row 7 - AdaBoostClassifier
$ws.Range("C7").Value = 0.4075

# Row 8 - LogisticRegression
$ws.Range("C8").Value = 0.8325

# Row 9 - SVC1
$ws.Range("C9").Value = 0.86

# Row 10 - SVC2
$ws.Range("C10").Value = 0.82

# Row 11 - SVC3
$ws.Range("C11").Value = 0.8575
